$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 20; $r++) {
    $ws.Range("I$r").Formula = '=C' + $r + '*C' + $r + '-$G$3*$G$3'
    $ws.Range("K$r").Formula = '=1/I' + $r
}

$ws.Range("K4:K20").Select()
